# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (column E) list for the worker is rebuilt in
# descending (most-recent-first) order instead of the previous
# ascending order. The "Valor Mora" (column F) values follow the rows,
# so the amounts that used to sit on the first/last period rows swap
# places as a side effect of the reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New descending period order for rows 16..51 (was ascending 1703..2003)
$periods = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706","1705","1703"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The "Valor Mora" amount that used to be on the first data row (1703 / 29509)
# now belongs to the last data row (now showing 1703), and the amount that
# used to be on the last data row (2003 / 48000) now belongs to the first
# data row (now showing 2003).
$ws.Range("F16").Value = 48000
$ws.Range("F51").Value = 29509
